# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-name suffixes to the concrete
#   format-version tags "_FV2304" / "_FV2310".
# - Turn the A1:U59 range into a real Excel Table ("Table1") with an
#   autofilter.
# - Freeze the header row (split/freeze pane under row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header row: "<name>_old" -> "<name>_FV2304", "<name>_new" -> "<name>_FV2310"
# ---------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2) Convert the used range into a native Excel table with an autofilter.
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------
# 3) Freeze the header row (pane split under row 1).
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
